$d = $word.ActiveDocument

# --- Split 1: "{m" -> "{" + "m" ----------------------------------------
# Locate the "{m" text (start of the field token) and split it into two
# separate runs at the boundary between "{" and "m", without altering the
# surrounding text or formatting.
$r1 = $d.Content
$r1.Find.Execute("{m", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$mChar = $d.Range($r1.Start + 1, $r1.Start + 2)
# Toggling a character formatting property on/off forces Word to break the
# run at this boundary while leaving the text and the rest of the run's
# formatting untouched.
$mChar.Bold = 1
$mChar.Bold = 0

# --- Split 2: ".fit(400, 400)}" -> ".fit(400, 400)" + "}" --------------
# Locate the closing ".fit(400, 400)}" text (end of the field token) and
# split it into two separate runs at the boundary right before the
# closing "}", preserving the orange character color applied to the run.
$r2 = $d.Content
$r2.Find.Execute(".fit(400, 400)}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$closeBrace = $d.Range($r2.End - 1, $r2.End)
$closeBrace.Bold = 1
$closeBrace.Bold = 0
